$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    # Force the cell to Text format so Excel does not coerce a
    # numeric-looking (or otherwise ambiguous) string into a number,
    # then clear the format back off so no stray style sticks around.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).ClearFormats()
}

Set-TextCell $ws "D2" "65.588.84"
$ws.Range("E2").Value = "  +0.19%  "

Set-TextCell $ws "D3" "2.660.87"
$ws.Range("E3").Value = "  +0.49%  "

$ws.Range("E4").Value = "  +0.05%  "

Set-TextCell $ws "D5" "600.74"
$ws.Range("E5").Value = "  -1.11%  "

Set-TextCell $ws "D6" "156.88"
$ws.Range("E6").Value = "  +0.48%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  +5.23%  "

$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("E10").Value = "  -0.68%  "

$ws.Range("E11").Value = "  -1.93%  "

$ws.Range("E12").Value = "  -0.24%  "

Set-TextCell $ws "D13" "29.32"
$ws.Range("E13").Value = "  -2.09%  "

Set-TextCell $ws "D14" "0.0000195"
$ws.Range("E14").Value = "  -4.32%  "

Set-TextCell $ws "D15" "3.136.90"
$ws.Range("E15").Value = "  +0.41%  "

Set-TextCell $ws "D16" "65.425.82"
$ws.Range("E16").Value = "  +0.21%  "

Set-TextCell $ws "D17" "2.670.76"
$ws.Range("E17").Value = "  +0.87%  "

Set-TextCell $ws "D18" "12.65"
$ws.Range("E18").Value = "  +0.04%  "

Set-TextCell $ws "D19" "4.80"
$ws.Range("E19").Value = "  -1.95%  "

$ws.Range("E20").Value = "  +2.23%  "

Set-TextCell $ws "D21" "351.17"
$ws.Range("E21").Value = "  -1.95%  "

$ws.Range("E22").Value = "  -0.12%  "

Set-TextCell $ws "D23" "69.34"
$ws.Range("E23").Value = "  -1.06%  "

Set-TextCell $ws "D24" "0.0000109"
$ws.Range("E24").Value = "  +3.66%  "

Set-TextCell $ws "D25" "9.72"
$ws.Range("E25").Value = "  +2.06%  "

$ws.Range("E26").Value = "  -4.14%  "

$ws.Range("E27").Value = "  -1.43%  "

Set-TextCell $ws "D28" "1.58"
$ws.Range("E28").Value = "  -3.09%  "

Set-TextCell $ws "D29" "8.06"
$ws.Range("E29").Value = "  -0.42%  "

$ws.Range("E30").Value = "  +0.61%  "

Set-TextCell $ws "D31" "532.62"
$ws.Range("E31").Value = "  -0.38%  "

$ws.Range("E32").Value = "  -3.14%  "

$ws.Range("E33").Value = "  -2.07%  "

Set-TextCell $ws "D34" "6.51"
$ws.Range("E34").Value = "  +2.14%  "

$ws.Range("E35").Value = "  -0.44%  "

Set-TextCell $ws "D36" "0.422"
$ws.Range("E36").Value = "  -2.10%  "

$ws.Range("E37").Value = "  -1.30%  "

$ws.Range("E38").Value = "  +0.05%  "

Set-TextCell $ws "D39" "159.15"
$ws.Range("E39").Value = "  -2.16%  "

$ws.Range("E40").Value = "  -3.37%  "

Set-TextCell $ws "D41" "0.999"
$ws.Range("E41").Value = "  -0.01%  "

Set-TextCell $ws "D42" "42.65"
$ws.Range("E42").Value = "  +1.56%  "

Set-TextCell $ws "D43" "164.36"
$ws.Range("E43").Value = "  -2.47%  "

$ws.Range("E44").Value = "  -1.88%  "

Set-TextCell $ws "D45" "0.0608"
$ws.Range("E45").Value = "  -0.25%  "

$ws.Range("E46").Value = "  -1.67%  "

Set-TextCell $ws "D47" "22.90"
$ws.Range("E47").Value = "  -0.37%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell $ws "D48" "0.641"
$ws.Range("E48").Value = "  -1.98%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws "D49" "0.0258"
$ws.Range("E49").Value = "  -2.15%  "

$ws.Range("E50").Value = "  +3.52%  "

Set-TextCell $ws "D51" "20.16"
$ws.Range("E51").Value = "  +2.41%  "

